# Add new field definitions ("Acici_Tipi", "Cift_Kafa_Var_Yok", "Cift_Kafa_Tipi")
# to the "fields" sheet, and their corresponding option lists to the
# "options" sheet, then leave the workbook's view state (active sheet /
# selections) the way the author left it when they were done editing.
#
# The cell-by-cell order below mirrors how the sheets were actually typed
# in (fields row 3 first, then the options rows for "Acici_Tipi_opts" filled
# in column-by-column, then fields rows 4-5, then the remaining options rows
# column-by-column) so that new entries land in the shared-string table in
# the same sequence as the original edit.

$wb = $excel.ActiveWorkbook
$fields = $wb.Worksheets.Item("fields")
$options = $wb.Worksheets.Item("options")

# ---------------------------------------------------------------------------
# 1. "fields" row 3 - Acici_Tipi
# ---------------------------------------------------------------------------
$fields.Range("A3").Value = "Acici_Tipi"
$fields.Range("B3").Value = "Acici_Tipi"
$fields.Range("C3").Value = "Açıcı Tipi"
$fields.Range("D3").Value = "select"
$fields.Range("E3").Value = $true
$fields.Range("F3").Value = "Acici_Tipi_opts"

# ---------------------------------------------------------------------------
# 2. "options" rows 7-8 - values for Acici_Tipi_opts (filled column-by-column)
# ---------------------------------------------------------------------------
$options.Range("A7").Value = "Acici_Tipi_opts"
$options.Range("A8").Value = "Acici_Tipi_opts"
$options.Range("B7").Value = "CMC"
$options.Range("B8").Value = "CMC(H)"
$options.Range("C7").Value = "Mekanik Tambur"
$options.Range("C8").Value = "Mekanikten Hidroliğe Çevirme"
$options.Range("D7").Value = 1
$options.Range("D8").Value = 2

# ---------------------------------------------------------------------------
# 3. "fields" rows 4-5 - Cift_Kafa_Var_Yok / Cift_Kafa_Tipi
# ---------------------------------------------------------------------------
$fields.Range("A4").Value = "Cift_Kafa"
$fields.Range("B4").Value = "Cift_Kafa_Var_Yok"
$fields.Range("C4").Value = "Mevcut mu?"
$fields.Range("D4").Value = "select"
$fields.Range("E4").Value = $true
$fields.Range("F4").Value = "Cift_Kafa_Var_Yok_opts"

$fields.Range("A5").Value = "Cift_Kafa"
$fields.Range("B5").Value = "Cift_Kafa_Tipi"
$fields.Range("C5").Value = "Çift Kafa Dönüş Tipi"
$fields.Range("D5").Value = "select"
$fields.Range("E5").Value = $true
$fields.Range("F5").Value = "Cift_Kafa_Tipi_opts"

# Column F now holds the longest text in the sheet ("Cift_Kafa_Var_Yok_opts") -
# widen it (best-fit) so the key names aren't truncated, matching the other
# "best fit" columns used throughout this workbook.
$fields.Columns.Item(6).AutoFit()
$fields.Columns.Item(6).ColumnWidth = 21.67

# ---------------------------------------------------------------------------
# 4. "options" rows 9-12 - values for Cift_Kafa_Var_Yok / Cift_Kafa_Tipi
#    (again filled column-by-column)
# ---------------------------------------------------------------------------
$options.Range("A9").Value = "Cift_Kafa_Var_Yok"
$options.Range("A10").Value = "Cift_Kafa_Var_Yok"
$options.Range("A11").Value = "Cift_Kafa_Tipi"
$options.Range("A12").Value = "Cift_Kafa_Tipi"
$options.Range("B9").Value = "Var"
$options.Range("B10").Value = "Yok"
$options.Range("B11").Value = "MT"
$options.Range("B12").Value = "AT"
$options.Range("C9").Value = "Çift Kafalı Dönüş"
$options.Range("C10").Value = "Tek Tambur"
$options.Range("C11").Value = "Manuel Döndürme"
$options.Range("C12").Value = "Otomatik Dönüş Sistemi"
$options.Range("D9").Value = 1
$options.Range("D10").Value = 2
$options.Range("D11").Value = 1
$options.Range("D12").Value = 2

# ---------------------------------------------------------------------------
# 5. Restore the cursor/selection on each sheet, then leave "fields" as the
#    active tab (mirrors the navigation the author did while editing).
# ---------------------------------------------------------------------------
$sections = $wb.Worksheets.Item("sections")
$sections.Range("D3").Select()

$options.Range("C16").Select()

$fields.Activate()
$fields.Range("J4").Select()
